$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.6132958801498127
$ws1.Range("C2").Value = 0.5641569459172853
$ws1.Range("D2").Value = 0.9962546816479401
$ws1.Range("E2").Value = 0.7203791469194313
$ws1.Range("F2").Value = 0.8639168561221175
$ws1.Range("G2").Value = 0.9677464493108515
$ws1.Range("H2").Value = 0.7977668363983222
$ws1.Range("I2").Value = 532
$ws1.Range("J2").Value = 411
$ws1.Range("K2").Value = 123
$ws1.Range("L2").Value = 2

# --- Sheet 2: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")

# Row 2: "0"
$ws2.Range("B2").Value = 0.984
$ws2.Range("C2").Value = 0.2303370786516854
$ws2.Range("D2").Value = 0.3732928679817906

# Row 3: "1"
$ws2.Range("B3").Value = 0.5641569459172853
$ws2.Range("C3").Value = 0.9962546816479401
$ws2.Range("D3").Value = 0.7203791469194313

# Row 4: "accuracy"
$ws2.Range("B4").Value = 0.6132958801498127
$ws2.Range("C4").Value = 0.6132958801498127
$ws2.Range("D4").Value = 0.6132958801498127
$ws2.Range("E4").Value = 0.6132958801498127

# Row 5: "macro avg"
$ws2.Range("B5").Value = 0.7740784729586426
$ws2.Range("C5").Value = 0.6132958801498127
$ws2.Range("D5").Value = 0.5468360074506109

# Row 6: "weighted avg"
$ws2.Range("B6").Value = 0.7740784729586426
$ws2.Range("C6").Value = 0.6132958801498127
$ws2.Range("D6").Value = 0.5468360074506109

# --- Sheet 3: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

# Row 2: "Actual 0"
$ws3.Range("B2").Value = 123
$ws3.Range("C2").Value = 411

# Row 3: "Actual 1"
$ws3.Range("B3").Value = 2
$ws3.Range("C3").Value = 532
